$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first anal")

# ---- Row 4 ----
$ws.Range("D4").Value = "7161,9026"
$ws.Range("E4").Value = 0.1249985694885254
$ws.Range("F4").Value = 0.07812595367431641
$ws.Range("G4").Value = 0.03125190734863281
$ws.Range("H4").Value = 975.3
$ws.Range("J4").Value = 991.5
$ws.Range("K4").Value = 2577.4
$ws.Range("N4").Value = 2612
$ws.Range("P4").Value = 646
$ws.Range("S4").Value = 940.5
$ws.Range("U4").Value = 991.5
$ws.Range("V4").Value = 2596.7
$ws.Range("X4").Value = 645
$ws.Range("Y4").Value = 2627.5
$ws.Range("AA4").Value = 639
$ws.Range("AE4").Value = 885
$ws.Range("AG4").Value = 991.5
$ws.Range("AH4").Value = 2613.2
$ws.Range("AJ4").Value = 661
$ws.Range("AK4").Value = 2665
$ws.Range("AM4").Value = 659

# ---- Row 5 ----
$ws.Range("D5").Value = "5347,8967"
$ws.Range("E5").Value = 0.140622615814209
$ws.Range("F5").Value = 0.09375572204589844
$ws.Range("G5").Value = 0.03124737739562988
$ws.Range("H5").Value = 943.2
$ws.Range("K5").Value = 2542.8
$ws.Range("M5").Value = 652
$ws.Range("N5").Value = 2530.4
$ws.Range("P5").Value = 653
$ws.Range("S5").Value = 901.5
$ws.Range("V5").Value = 2559.3
$ws.Range("X5").Value = 644
$ws.Range("Y5").Value = 2556.9
$ws.Range("AA5").Value = 644
$ws.Range("AE5").Value = 829.5
$ws.Range("AH5").Value = 2601.3
$ws.Range("AJ5").Value = 666
$ws.Range("AK5").Value = 2586.9
$ws.Range("AM5").Value = 662
